$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the workbook window (best-effort; reflects the saved bookView size).
try {
    $win = $excel.ActiveWindow
    $win.Width = 18468
    $win.Height = 9048
} catch {
    # ignore if not supported by the host
}

# Move the active cell selection on the sheet.
$ws.Range("I8").Select()

# Fill in the local-explanation values (B:G) for rows 3-6.
$data = @{
    3 = @(15, 2, 3, 2, 3, 4)
    4 = @(16, 3, 3, 3, 3, 4)
    5 = @(16, 2, 3, 2, 3, 4)
    6 = @(17, 3, 3, 3, 3, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $cols = @("B", "C", "D", "E", "F", "G")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
